{"js": "// Replace each arithmetic expression in the worksheet table with its\n// updated value, in row-major order (20 rows x 5 columns), matching the\n// target revision. Formatting (fonts, size, run/paragraph properties) is\n// left untouched because we only rewrite the table's text values.\nconst table = context.document.body.tables.getFirst();\ntable.load(\"values\");\nawait context.sync();\n\nconst newValues = [\n  [\"16+21=\", \"0+33=\", \"72-7=\", \"74-62=\", \"65-27=\"],\n  [\"2+77=\", \"52+36=\", \"4+85=\", \"20+57=\", \"51-17=\"],\n  [\"49-49=\", \"42+29=\", \"75-9=\", \"60-3=\", \"12+36=\"],\n  [\"37+57=\", \"92-46=\", \"37-4=\", \"67-59=\", \"34-29=\"],\n  [\"39+5=\", \"47-5=\", \"38+56=\", \"61-37=\", \"10+0=\"],\n  [\"10+57=\", \"90-10=\", \"31-10=\", \"23+8=\", \"8+60=\"],\n  [\"26+29=\", \"45-17=\", \"80+11=\", \"7+9=\", \"12+79=\"],\n  [\"37+34=\", \"44+17=\", \"39+24=\", \"29-21=\", \"26+73=\"],\n  [\"50-18=\", \"68+9=\", \"43+20=\", \"83-13=\", \"76-26=\"],\n  [\"54-5=\", \"52-36=\", \"70-2=\", \"15+81=\", \"59-15=\"],\n  [\"98-60=\", \"30+68=\", \"1+57=\", \"64+28=\", \"14+12=\"],\n  [\"31+17=\", \"99-59=\", \"41+43=\", \"11+54=\", \"73+11=\"],\n  [\"59+10=\", \"84-34=\", \"7+90=\", \"14+1=\", \"11+31=\"],\n  [\"21+13=\", \"10+26=\", \"27+41=\", \"2+46=\", \"5+80=\"],\n  [\"34+10=\", \"75-48=\", \"54+44=\", \"89-35=\", \"32+44=\"],\n  [\"43+32=\", \"85-52=\", \"6+46=\", \"7+4=\", \"89-54=\"],\n  [\"34-4=\", \"25+24=\", \"44+39=\", \"3+82=\", \"15+55=\"],\n  [\"26+30=\", \"10+57=\", \"77-59=\", \"40+50=\", \"2+2=\"],\n  [\"5+51=\", \"78-24=\", \"41+8=\", \"86+1=\", \"18-9=\"],\n  [\"1+42=\", \"93-68=\", \"74+19=\", \"51+21=\", \"47+10=\"]\n];\n\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# Replace each arithmetic expression in the worksheet table with its\n# updated value, in row-major order (20 rows x 5 columns), matching the\n# target revision. Setting Cell.Range.Text rewrites only the cell's text\n# run, so formatting (fonts, size, run/paragraph properties) is preserved.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$newValues = @(\n    @(\"16+21=\", \"0+33=\", \"72-7=\", \"74-62=\", \"65-27=\"),\n    @(\"2+77=\", \"52+36=\", \"4+85=\", \"20+57=\", \"51-17=\"),\n    @(\"49-49=\", \"42+29=\", \"75-9=\", \"60-3=\", \"12+36=\"),\n    @(\"37+57=\", \"92-46=\", \"37-4=\", \"67-59=\", \"34-29=\"),\n    @(\"39+5=\", \"47-5=\", \"38+56=\", \"61-37=\", \"10+0=\"),\n    @(\"10+57=\", \"90-10=\", \"31-10=\", \"23+8=\", \"8+60=\"),\n    @(\"26+29=\", \"45-17=\", \"80+11=\", \"7+9=\", \"12+79=\"),\n    @(\"37+34=\", \"44+17=\", \"39+24=\", \"29-21=\", \"26+73=\"),\n    @(\"50-18=\", \"68+9=\", \"43+20=\", \"83-13=\", \"76-26=\"),\n    @(\"54-5=\", \"52-36=\", \"70-2=\", \"15+81=\", \"59-15=\"),\n    @(\"98-60=\", \"30+68=\", \"1+57=\", \"64+28=\", \"14+12=\"),\n    @(\"31+17=\", \"99-59=\", \"41+43=\", \"11+54=\", \"73+11=\"),\n    @(\"59+10=\", \"84-34=\", \"7+90=\", \"14+1=\", \"11+31=\"),\n    @(\"21+13=\", \"10+26=\", \"27+41=\", \"2+46=\", \"5+80=\"),\n    @(\"34+10=\", \"75-48=\", \"54+44=\", \"89-35=\", \"32+44=\"),\n    @(\"43+32=\", \"85-52=\", \"6+46=\", \"7+4=\", \"89-54=\"),\n    @(\"34-4=\", \"25+24=\", \"44+39=\", \"3+82=\", \"15+55=\"),\n    @(\"26+30=\", \"10+57=\", \"77-59=\", \"40+50=\", \"2+2=\"),\n    @(\"5+51=\", \"78-24=\", \"41+8=\", \"86+1=\", \"18-9=\"),\n    @(\"1+42=\", \"93-68=\", \"74+19=\", \"51+21=\", \"47+10=\")\n)\n\nfor ($r = 1; $r -le $newValues.Count; $r++) {\n    $row = $newValues[$r - 1]\n    for ($c = 1; $c -le $row.Count; $c++) {\n        $t.Cell($r, $c).Range.Text = $row[$c - 1]\n    }\n}\n"}
